$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 284.5
$ws.Range("I41").Value = 130
$ws.Range("J41").Value = 439
$ws.Range("K41").Value = 130
$ws.Range("L41").Value = 439
$ws.Range("M41").Value = 310
$ws.Range("N41").Value = -1319
$ws.Range("H62").Value = 7353.636
$ws.Range("J62").Value = 7861.375
$ws.Range("L62").Value = 7861.375
$ws.Range("N62").Value = -9109.375
$ws.Range("H65").Value = 7353.636
$ws.Range("J65").Value = 7861.375
$ws.Range("L65").Value = 39306.875
$ws.Range("N65").Value = -45546.875
$ws.Range("H98").Value = 1666.6666
$ws.Range("I98").Value = 750
$ws.Range("J98").Value = 3500
$ws.Range("K98").Value = 750
$ws.Range("L98").Value = 3500
$ws.Range("M98").Value = 748
$ws.Range("N98").Value = -6496
$ws.Range("H111").Value = 3219.2666
$ws.Range("I111").Value = 4132.1113
$ws.Range("J111").Value = 1850
$ws.Range("K111").Value = 12396.3339
$ws.Range("L111").Value = 5550
$ws.Range("M111").Value = -9329.333899999998
$ws.Range("N111").Value = -11684
$ws.Range("H122").Value = 1666.6666
$ws.Range("I122").Value = 750
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 2250
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = 200
$ws.Range("N122").Value = -15400
$ws.Range("H127").Value = 1419.1818
$ws.Range("I127").Value = 999.5
$ws.Range("J127").Value = 1922.8
$ws.Range("K127").Value = 2998.5
$ws.Range("L127").Value = 5768.4
$ws.Range("M127").Value = 1961.5
$ws.Range("N127").Value = -15688.4
$ws.Range("H129").Value = 1133.317
$ws.Range("J129").Value = 1265.3143
$ws.Range("L129").Value = 3795.9429
$ws.Range("N129").Value = -13795.9429
$ws.Range("H132").Value = 3327.9707
$ws.Range("I132").Value = 3255.2856
$ws.Range("K132").Value = 9765.856800000001
$ws.Range("M132").Value = -7235.856800000001
$ws.Range("H135").Value = 19233020
$ws.Range("I135").Value = 641.381
$ws.Range("K135").Value = 5772.429
$ws.Range("M135").Value = -3237.429

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3065.8125
$ws.Range("I45").Value = 2776.375
$ws.Range("J45").Value = 3355.25
$ws.Range("K45").Value = 2776.375
$ws.Range("L45").Value = 3355.25
$ws.Range("M45").Value = -2399.375
$ws.Range("N45").Value = -4109.25
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H82").Value = 21500
$ws.Range("I82").Value = 15000
$ws.Range("J82").Value = 28000
$ws.Range("K82").Value = 15000
$ws.Range("L82").Value = 28000
$ws.Range("M82").Value = -14639
$ws.Range("N82").Value = -28722
$ws.Range("H85").Value = 21500
$ws.Range("I85").Value = 15000
$ws.Range("J85").Value = 28000
$ws.Range("K85").Value = 15000
$ws.Range("L85").Value = 28000
$ws.Range("M85").Value = -13752
$ws.Range("N85").Value = -30496
$ws.Range("H122").Value = 2851.0417
$ws.Range("J122").Value = 4833
$ws.Range("L122").Value = 14499
$ws.Range("N122").Value = -19399
$ws.Range("H132").Value = 21262.814
$ws.Range("I132").Value = 2507.5264
$ws.Range("K132").Value = 7522.5792
$ws.Range("M132").Value = -4992.5792

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 2079.9
$ws.Range("I11").Value = 2079.9
$ws.Range("K11").Value = 2079.9
$ws.Range("M11").Value = -1939.9
$ws.Range("H94").Value = 914.2105
$ws.Range("I94").Value = 735.38464
$ws.Range("K94").Value = 735.38464
$ws.Range("M94").Value = -284.38464
$ws.Range("H134").Value = 4072.2307
$ws.Range("I134").Value = 4781.619
$ws.Range("J134").Value = 1092.8
$ws.Range("K134").Value = 14344.857
$ws.Range("L134").Value = 3278.4
$ws.Range("M134").Value = -11809.857
$ws.Range("N134").Value = -8348.4

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1023.1667
$ws.Range("I16").Value = 812.8570999999999
$ws.Range("K16").Value = 812.8570999999999
$ws.Range("M16").Value = -525.8570999999999
$ws.Range("H22").Value = 554.3333
$ws.Range("J22").Value = 751
$ws.Range("L22").Value = 751
$ws.Range("N22").Value = -1451
$ws.Range("H31").Value = 2989.9119
$ws.Range("I31").Value = 1250.75
$ws.Range("J31").Value = 3525.0386
$ws.Range("K31").Value = 1250.75
$ws.Range("L31").Value = 3525.0386
$ws.Range("M31").Value = -955.75
$ws.Range("N31").Value = -4115.0386
$ws.Range("H34").Value = 2989.9119
$ws.Range("I34").Value = 1250.75
$ws.Range("J34").Value = 3525.0386
$ws.Range("K34").Value = 1250.75
$ws.Range("L34").Value = 3525.0386
$ws.Range("M34").Value = -1048.75
$ws.Range("N34").Value = -3929.0386
$ws.Range("H86").Value = 19948.875
$ws.Range("I86").Value = 3966.6667
$ws.Range("J86").Value = 29538.2
$ws.Range("K86").Value = 3966.6667
$ws.Range("L86").Value = 29538.2
$ws.Range("M86").Value = -2843.6667
$ws.Range("N86").Value = -31784.2
$ws.Range("H89").Value = 19948.875
$ws.Range("I89").Value = 3966.6667
$ws.Range("J89").Value = 29538.2
$ws.Range("K89").Value = 19833.3335
$ws.Range("L89").Value = 147691
$ws.Range("M89").Value = -14217.3335
$ws.Range("N89").Value = -158923
$ws.Range("H107").Value = 1537.5652
$ws.Range("I107").Value = 1361.25
$ws.Range("J107").Value = 1631.6
$ws.Range("K107").Value = 1361.25
$ws.Range("L107").Value = 1631.6
$ws.Range("M107").Value = 558.75
$ws.Range("N107").Value = -5471.6
$ws.Range("H113").Value = 1023.1667
$ws.Range("I113").Value = 812.8570999999999
$ws.Range("K113").Value = 812.8570999999999
$ws.Range("M113").Value = 1357.1429

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 1106
$ws.Range("I17").Value = 190
$ws.Range("J17").Value = 1335
$ws.Range("K17").Value = 570
$ws.Range("L17").Value = 4005
$ws.Range("M17").Value = -401
$ws.Range("N17").Value = -4343
$ws.Range("H40").Value = 140
$ws.Range("I40").Value = 76
$ws.Range("K40").Value = 304
$ws.Range("M40").Value = -235
$ws.Range("H80").Value = 8425
$ws.Range("J80").Value = 8842.385
$ws.Range("L80").Value = 26527.155
$ws.Range("N80").Value = -28399.155
$ws.Range("H83").Value = 8425
$ws.Range("J83").Value = 8842.385
$ws.Range("L83").Value = 79581.465
$ws.Range("N83").Value = -88941.465
$ws.Range("H121").Value = 936.03125
$ws.Range("I121").Value = 489.125
$ws.Range("J121").Value = 1085
$ws.Range("K121").Value = 1467.375
$ws.Range("L121").Value = 3255
$ws.Range("M121").Value = -157.375
$ws.Range("N121").Value = -5875
$ws.Range("H131").Value = 704.3737
$ws.Range("J131").Value = 721.10754
$ws.Range("L131").Value = 2163.32262
$ws.Range("N131").Value = -12243.32262
$ws.Range("H140").Value = 2596.8
$ws.Range("I140").Value = 817.7778
$ws.Range("J140").Value = 4052.3635
$ws.Range("K140").Value = 2453.3334
$ws.Range("L140").Value = 12157.0905
$ws.Range("M140").Value = 2726.6666
$ws.Range("N140").Value = -22517.0905

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1974.76
$ws.Range("I113").Value = 1623.7333
$ws.Range("J113").Value = 2501.3
$ws.Range("K113").Value = 1623.7333
$ws.Range("L113").Value = 2501.3
$ws.Range("M113").Value = 546.2666999999999
$ws.Range("N113").Value = -6841.3

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5087.75
$ws.Range("I22").Value = 3450.6667
$ws.Range("K22").Value = 3450.6667
$ws.Range("M22").Value = -3155.6667
$ws.Range("H27").Value = 5087.75
$ws.Range("I27").Value = 3450.6667
$ws.Range("K27").Value = 3450.6667
$ws.Range("M27").Value = -3343.6667
$ws.Range("H63").Value = 18000
$ws.Range("J63").Value = 18000
$ws.Range("L63").Value = 18000
$ws.Range("N63").Value = -19498
$ws.Range("H66").Value = 18000
$ws.Range("J66").Value = 18000
$ws.Range("L66").Value = 54000
$ws.Range("N66").Value = -61488
$ws.Range("H122").Value = 1403688.8
$ws.Range("I122").Value = 1963344.2
$ws.Range("J122").Value = 4549.75
$ws.Range("K122").Value = 5890032.6
$ws.Range("L122").Value = 13649.25
$ws.Range("M122").Value = -5887582.6
$ws.Range("N122").Value = -18549.25
$ws.Range("H132").Value = 863080.5600000001
$ws.Range("I132").Value = 1206713.5
$ws.Range("J132").Value = 3998.25
$ws.Range("K132").Value = 3620140.5
$ws.Range("L132").Value = 11994.75
$ws.Range("M132").Value = -3617610.5
$ws.Range("N132").Value = -17054.75

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2013.7693
$ws.Range("I132").Value = 1210.125
$ws.Range("J132").Value = 3299.6
$ws.Range("K132").Value = 3630.375
$ws.Range("L132").Value = 9898.799999999999
$ws.Range("M132").Value = -1100.375
$ws.Range("N132").Value = -14958.8
